$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Name", "Age", "Gender"),
    @("Josh", 23, "Male"),
    @("Zach", 20, "Male"),
    @("Jesse", 17, "Male"),
    @("Andrew", 14, "Male"),
    @("Chris", 11, "Male")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

$ws.Range("A7").Select()
